# Update cryptos list (price + volume changes scraped on Sat Jun 24 04:57:05 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "30.728.76" / "0.000007990" - force Text format
# so the COM layer does not reinterpret these as numbers (stripping the
# thousands-dot grouping or trailing zeros).
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.728.76'
$ws.Range('E2').Value = '  +2.53%  '
$ws.Range('D3').Value = '1.893.70'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '246.69'
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').Value = '0.4928'
$ws.Range('E7').Value = '  -1.15%  '
$ws.Range('D8').Value = '0.2957'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').Value = '0.06807'
$ws.Range('E9').Value = '  +2.67%  '
$ws.Range('D10').Value = '1.893.37'
$ws.Range('D11').Value = '17.30'
$ws.Range('E11').Value = '  +3.37%  '
$ws.Range('D12').Value = '92.22'
$ws.Range('D13').Value = '0.07257'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').Value = '0.6825'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').Value = '5.095'
$ws.Range('E15').Value = '  +4.32%  '
$ws.Range('D16').Value = '30.711.40'
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('D17').Value = '0.000007990'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').Value = '13.31'
$ws.Range('E18').Value = '  +4.17%  '
$ws.Range('D20').Value = '2.139.96'
$ws.Range('D21').Value = '1.001'
$ws.Range('D22').Value = '4.849'
$ws.Range('E22').Value = '  +1.82%  '
$ws.Range('D23').Value = '191.35'
$ws.Range('E23').Value = '  +34.90%  '
$ws.Range('D24').Value = '6.062'
$ws.Range('E24').Value = '  +7.22%  '
$ws.Range('D25').Value = '9.402'
$ws.Range('E25').Value = '  +3.79%  '
$ws.Range('D26').Value = '155.90'
$ws.Range('E26').Value = '  +4.52%  '
$ws.Range('D27').Value = '19.15'
$ws.Range('E27').Value = '  +11.28%  '
$ws.Range('D28').Value = '1.928'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').Value = '1.400'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').Value = '4.379'
$ws.Range('E30').Value = '  +4.97%  '
$ws.Range('D31').Value = '0.09015'
$ws.Range('D32').Value = '4.045'
$ws.Range('E32').Value = '  +2.60%  '
$ws.Range('E33').Value = '  +2.64%  '
$ws.Range('D34').Value = '0.7467'
$ws.Range('E34').Value = '  +5.28%  '
$ws.Range('D35').Value = '1.126'
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('D36').Value = '2.724'
$ws.Range('E36').Value = '  +2.23%  '
$ws.Range('D37').Value = '0.01862'
$ws.Range('E37').Value = '  +6.05%  '
$ws.Range('D38').Value = '2.676'
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('D39').Value = '2.164'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').Value = '0.9386'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').Value = '0.4443'
$ws.Range('E41').Value = '  +4.49%  '
$ws.Range('D42').Value = '106.48'
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('D43').Value = '5.802'
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').Value = '7.697'
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('D46').Value = '0.1342'
$ws.Range('E46').Value = '  +6.67%  '
$ws.Range('D47').Value = '0.05856'
$ws.Range('E47').Value = '  +3.63%  '
$ws.Range('D48').Value = '8.751'
$ws.Range('E48').Value = '  +6.64%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.433'
$ws.Range('E49').Value = '  +7.48%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').Value = '0.3958'
$ws.Range('E50').Value = '  +5.25%  '
$ws.Range('D51').Value = '33.62'
$ws.Range('E51').Value = '  +3.68%  '
